$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that follows the title heading.
#    (It contained a bold "Meta description" run plus the description text.)
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Near the end of the document, insert a new bold paragraph containing the
#    title text right before the final paragraph, and change the final
#    paragraph's text from the DALLE image prompt to the meta-description
#    text (keeping its existing italic formatting).
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$targetRange = $lastPara.Range

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruity Mania Free: A Classic Land-Based Casino Experience</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the classic charm of Fruity Mania, an online slot game with straightforward gameplay, lively atmosphere, and decent payouts. Play now for free!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$targetRange.InsertXML($xmlFragment)

# InsertXML leaves one stray empty paragraph behind (the original paragraph
# mark that used to terminate the old last paragraph). Merge it away by
# deleting the mark that separates the second-to-last paragraph from it.
$newCount = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($newCount - 1)
$trailing = $d.Paragraphs.Item($newCount)
$cleanupRange = $d.Range($secondLast.Range.End - 1, $trailing.Range.End)
$cleanupRange.Delete()
